# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" (fund-level detail, same layout as the
#   "2021-Q3"/"2021-Q4" sheets) positioned right after "2021-Q4" and before
#   the "总计" summary sheet.
# - Add a new top row to "总计" for the 2022-Q1 roll-up, pushing the
#   existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as literal text (no numeric/leading
# -zero coercion), and make sure the cell is left with NO extra style
# applied (matches the plain, un-styled data cells used throughout this
# workbook).
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Match the page-setup/column-width defaults used by the rest of this
# workbook (0.75in/1in/0.5in margins, standard 8.43 column width) instead
# of the generic blank-sheet defaults a freshly Add()-ed sheet starts with.
# Must run before any Range writes on the sheet, otherwise the engine
# bakes in a (harmless but non-matching) explicit baseColWidth.
function Set-SheetDefaults {
    param($sheet)
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
    $sheet.StandardWidth = 8.43
}

# -----------------------------------------------------------------
# 1) Rebuild the sheet order/ids so the new sheet lines up exactly:
#    2021-Q3 (id1), 2021-Q4 (id2), 2022-Q1 (id3), 总计 (id4)
# -----------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# remember formatting source sheet before we touch anything
$fmtSource = $q4

# Drop the old "总计" sheet - it gets recreated below so it is last again
# (and its internal id becomes 4, leaving id 3 free for the new sheet).
$oldTotal.Delete()

# New fund-detail sheet, inserted right after "2021-Q4".
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"
Set-SheetDefaults $q1

# Recreate "总计" right after the new sheet, at the end again.
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"
Set-SheetDefaults $total

# -----------------------------------------------------------------
# 2) Populate "2022-Q1" (fund-level detail) using the same column
#    layout/formatting as "2021-Q4".
# -----------------------------------------------------------------
$fmtSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSource.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "009531"
$q1.Range("C2").Value = "九泰锐和18个月定期开放混合"
Set-TextValue $q1.Range("D2") "1.97"
Set-TextValue $q1.Range("E2") "70.75"
Set-TextValue $q1.Range("F2") "4.38"
Set-TextValue $q1.Range("G2") "0.0863"
$q1.Range("H2").Value = 4

$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "010764"
$q1.Range("C3").Value = "九泰锐升18个月封闭运作混合"
Set-TextValue $q1.Range("D3") "3.15"
Set-TextValue $q1.Range("E3") "78.81"
Set-TextValue $q1.Range("F3") "2.73"
Set-TextValue $q1.Range("G3") "0.0860"
$q1.Range("H3").Value = 5

# -----------------------------------------------------------------
# 3) Rebuild "总计": same 3 columns, now with an extra "2022-Q1" row on
#    top of the previous "2021-Q4" / "2021-Q3" rows.
# -----------------------------------------------------------------
$fmtSource.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtSource.Range("A2:A3").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.17

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.26

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.18

# Restore the originally active sheet/tab ("2021-Q3" was first/active
# before this edit, and none of the diff's bookViews changed).
$wb.Worksheets.Item("2021-Q3").Activate()

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Output ("  " + $s.Index.ToString() + ": " + $s.Name)
}
